$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Update the "onblock" value for the "Jumping Left Kick" move (C39) which
# drives the dependent formulas in C15/F15/F39 via shared formulas.
$ws.Range("C39").Value = 16

# Update the selected cell in the sheet view to reflect where the edit was made.
$ws.Range("C39").Select()
